$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1) Extend the "User Home:" paragraph with two more line-break runs ---
$pUserHome = $d.Paragraphs.Item(4)
$pUserHome.Range.InsertAfter([char]11)
$pUserHome.Range.InsertAfter([char]11 + "Should not increase the number of product by increasing the number of same product in the cart.")

# --- 2) Blank separator paragraph ---
$pUserHome.Range.InsertParagraphAfter() | Out-Null
$pBlank1 = $pUserHome.Next()
$pBlank1.Range.InsertXML("<w:p $wNs/>") | Out-Null

# --- 3) "Fix the admin routes..." paragraph (with proofErr around "textfield") ---
$pBlank1.Range.InsertParagraphAfter() | Out-Null
$pFix = $pBlank1.Next()
$fixXml = "<w:p $wNs>" +
  "<w:r><w:t>Fix the admin routes so that user can" + [char]0x2019 + "t access the admin links.</w:t></w:r>" +
  "<w:r><w:br/></w:r>" +
  "<w:r><w:br/><w:t xml:space=`"preserve`">Fix : the sub categories </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>textfield</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:t xml:space=`"preserve`"> of categories looses focus after inserting single letter.</w:t></w:r>" +
  "</w:p>"
$pFix.Range.InsertXML($fixXml) | Out-Null

# --- 4) Blank separator paragraph ---
$pFix.Range.InsertParagraphAfter() | Out-Null
$pBlank2 = $pFix.Next()
$pBlank2.Range.InsertXML("<w:p $wNs/>") | Out-Null

# --- 5) "Home.jsx is loading continuously." paragraph (proofErr around "Home.jsx") ---
$pBlank2.Range.InsertParagraphAfter() | Out-Null
$pHome = $pBlank2.Next()
$homeXml = "<w:p $wNs>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>Home.jsx</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:t xml:space=`"preserve`"> is loading continuously. </w:t></w:r>" +
  "</w:p>"
$pHome.Range.InsertXML($homeXml) | Out-Null

Write-Host "done"
